# Daily attendance processing - 2025-12-05 05:27:10
# Normalize the "Recorded By" (column G) list ordering on the
# "Session Analysis Results" sheet. Each cell holds a comma-separated list
# of recorder identifiers (names / emails); re-order the items within each
# list to the canonical ordering below (values themselves are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# before-order -> after-order for the distinct "Recorded By" combinations
# that need to be re-sequenced.
$map = @{
    "system, System, backup@backdoor.com" = "System, backup@backdoor.com, system";
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Value()
    if ($null -eq $current) {
        continue
    }
    $key = [string]$current
    if ($map.ContainsKey($key)) {
        $cell.Value = $map[$key]
    }
}
